$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text formatting so numeric-looking
# strings (e.g. "546.69", "2.365.30") are not auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("B30").NumberFormat = "@"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = '60.961.36'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '2.365.30'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '546.69'
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").Value = '137.20'
$ws.Range("E6").Value = '  -3.65%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.524'
$ws.Range("E8").Value = '  -2.96%  '
$ws.Range("D9").Value = '2.364.47'
$ws.Range("E9").Value = '  -1.21%  '
$ws.Range("D10").Value = '0.106'
$ws.Range("E10").Value = '  +0.04%  '
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("D12").Value = '5.34'
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("D13").Value = '0.346'
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").Value = '24.90'
$ws.Range("E14").Value = '  -2.08%  '
$ws.Range("D15").Value = '2.775.30'
$ws.Range("E15").Value = '  -1.63%  '
$ws.Range("D16").Value = '0.0000165'
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("D17").Value = '60.890.17'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '2.365.52'
$ws.Range("E18").Value = '  -0.97%  '
$ws.Range("D19").Value = '10.72'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '319.41'
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = '4.12'
$ws.Range("E21").Value = '  -0.88%  '
$ws.Range("D22").Value = '6.63'
$ws.Range("E22").Value = '  -2.04%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = '64.03'
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("D25").Value = '1.64'
$ws.Range("E25").Value = '  -16.19%  '
$ws.Range("D26").Value = '8.22'
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").Value = '2.472.80'
$ws.Range("E28").Value = '  -1.35%  '
$ws.Range("D29").Value = '8.08'
$ws.Range("E29").Value = '  -0.49%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = '0.149'
$ws.Range("E30").Value = '  +1.62%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '0.0₃0876'
$ws.Range("E31").Value = '  -6.80%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = '1.37'
$ws.Range("E32").Value = '  -5.29%  '
$ws.Range("D33").Value = '496.73'
$ws.Range("E33").Value = '  -7.90%  '
$ws.Range("E34").Value = '  -1.41%  '
$ws.Range("D35").Value = '1.50'
$ws.Range("E35").Value = '  -5.82%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").Value = '4.66'
$ws.Range("E37").Value = '  -1.83%  '
$ws.Range("D38").Value = '1.88'
$ws.Range("E38").Value = '  +2.91%  '
$ws.Range("D39").Value = '0.377'
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("B40").Value = 'EthereumClassic'
$ws.Range("C40").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D40").Value = '18.47'
$ws.Range("E40").Value = '  +1.76%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").Value = '5.31'
$ws.Range("E41").Value = '  -5.19%  '
$ws.Range("D42").Value = '145.13'
$ws.Range("E42").Value = '  +4.40%  '
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").Value = '41.07'
$ws.Range("E44").Value = '  +1.80%  '
$ws.Range("D45").Value = '145.65'
$ws.Range("E45").Value = '  +2.50%  '
$ws.Range("D46").Value = '3.57'
$ws.Range("E46").Value = '  -1.36%  '
$ws.Range("D47").Value = '2.03'
$ws.Range("E47").Value = '  -8.85%  '
$ws.Range("D48").Value = '0.0520'
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("D49").Value = '19.15'
$ws.Range("E49").Value = '  -5.66%  '
$ws.Range("D50").Value = '0.572'
$ws.Range("E50").Value = '  -0.98%  '
$ws.Range("D51").Value = '0.0905'
$ws.Range("E51").Value = '  -0.33%  '
